{"js": "// Insert template placeholders for the spelled-out peso amounts, and\n// replace the trailing tab with the \"Oaxaca de Juarez, con fecha ${fecha_completa}\"\n// text, in the contract body (word/document.xml).\n\nconst body = context.document.body;\n\n// Helper: find a unique text fragment in the body and replace it in place\n// (InsertLocation.Replace keeps the formatting of the run(s) it overlaps).\nasync function replaceOnce(searchText, newText) {\n  const results = body.search(searchText, { matchCase: true, matchWildcards: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) \" (Quince mil pesos 00/100 M.N. ) con un I.V.A de $ \"\n//    -> \" (${cantidad_letras} pesos 00/100 M.N. ) con un I.V.A de $ \"\nawait replaceOnce(\"Quince mil \", \"${cantidad_letras} \");\n\n// 2) \" (Dos mil cuatrocientos pesos 00/100 M.N.) haciendo un  monto total de $\"\n//    -> \" (${iva_letras} pesos 00/100 M.N.) haciendo un  monto total de $\"\nawait replaceOnce(\"Dos mil cuatrocientos \", \"${iva_letras} \");\n\n// 3) \" (Diecisiete mil cuatrocientos pesos 00/100 M.N.), el 50% equivalente a la cantidad de $\"\n//    -> \" (${total_letras}  pesos 00/100 M.N.), el 50% equivalente a la cantidad de $\"\nawait replaceOnce(\"Diecisiete mil cuatrocientos \", \"${total_letras}  \");\n\n// 4) \" (Ocho mil setecientos pesos 00/100 M.N.) deber\u00e1 cubrirlo al momento ...\"\n//    -> \" (${total_mitad_letras} pesos 00/100 M.N.) deber\u00e1 cubrirlo al momento ...\"\nawait replaceOnce(\n  \"Ocho mil setecientos pesos 00/100 M.N.) deber\u00e1\",\n  \"${total_mitad_letras} pesos 00/100 M.N.) deber\u00e1\"\n);\n\n// 5) \" (Ocho mil setecientos pesos 00/100 M.N.) al t\u00e9rmino del servicio ...\"\n//    -> \" (${total_restante_letras} pesos 00/100 M.N.) al t\u00e9rmino del servicio ...\"\nawait replaceOnce(\n  \"Ocho mil setecientos pesos 00/100 M.N.) al t\u00e9rmino\",\n  \"${total_restante_letras} pesos 00/100 M.N.) al t\u00e9rmino\"\n);\n\n// 6) The lone tab character (\"... en la ciudad de [TAB].\") becomes\n//    \"... en la ciudad de Oaxaca de Juarez, con fecha ${fecha_completa}.\"\n//    Replacing just the tab keeps the single-underline formatting of that run.\nawait replaceOnce(\"\\t\", \"Oaxaca de Juarez, con fecha ${fecha_completa}\");\n", "ps1": "# Insert template placeholders for the spelled-out peso amounts, and\n# replace the trailing tab with \"Oaxaca de Juarez, con fecha ${fecha_completa}\"\n# in the contract body (word/document.xml).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($findText, $replaceText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue = 1, wdReplaceOne = 1 (we only ever expect one hit per call)\n    $find.Execute([ref]$findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1) | Out-Null\n}\n\n# 1) \" (Quince mil pesos 00/100 M.N. ) con un I.V.A de $ \"\n#    -> \" (${cantidad_letras} pesos 00/100 M.N. ) con un I.V.A de $ \"\nReplace-Once \"Quince mil \" \"`${cantidad_letras} \"\n\n# 2) \" (Dos mil cuatrocientos pesos 00/100 M.N.) haciendo un  monto total de $\"\n#    -> \" (${iva_letras} pesos 00/100 M.N.) haciendo un  monto total de $\"\nReplace-Once \"Dos mil cuatrocientos \" \"`${iva_letras} \"\n\n# 3) \" (Diecisiete mil cuatrocientos pesos 00/100 M.N.), el 50% equivalente a la cantidad de $\"\n#    -> \" (${total_letras}  pesos 00/100 M.N.), el 50% equivalente a la cantidad de $\"\nReplace-Once \"Diecisiete mil cuatrocientos \" \"`${total_letras}  \"\n\n# 4) \" (Ocho mil setecientos pesos 00/100 M.N.) deber\u00e1 cubrirlo al momento ...\"\n#    -> \" (${total_mitad_letras} pesos 00/100 M.N.) deber\u00e1 cubrirlo al momento ...\"\nReplace-Once \"Ocho mil setecientos pesos 00/100 M.N.) deber\" \"`${total_mitad_letras} pesos 00/100 M.N.) deber\"\n\n# 5) \" (Ocho mil setecientos pesos 00/100 M.N.) al t\u00e9rmino del servicio ...\"\n#    -> \" (${total_restante_letras} pesos 00/100 M.N.) al t\u00e9rmino del servicio ...\"\nReplace-Once \"Ocho mil setecientos pesos 00/100 M.N.) al t\" \"`${total_restante_letras} pesos 00/100 M.N.) al t\"\n\n# 6) The lone tab character (\"... en la ciudad de [TAB].\") becomes\n#    \"... en la ciudad de Oaxaca de Juarez, con fecha ${fecha_completa}.\"\n#    \"^t\" is Word's wildcard-free Find code for a tab character; replacing just\n#    the tab keeps the single-underline formatting of that run.\nReplace-Once \"^t\" \"Oaxaca de Juarez, con fecha `${fecha_completa}\"\n"}
